# Applies the "Updated cryptos list" data refresh to Sheet1.
# For each changed row, the Price (D) and/or Volume(1h) (E) columns are
# updated; row 50 additionally gets a new coin (BabyDogeCoin) inserted
# ahead of Algorand, pushing the former row 50/51 entries down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.699.67"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").Value = "1.644.60"

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.18"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.533"
$ws.Range("E6").Value = "  +3.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.10"
$ws.Range("E8").Value = "  -2.12%  "

$ws.Range("E9").Value = "  -0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0612"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("E11").Value = "  +1.66%  "

$ws.Range("D12").Value = "1.877.13"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").Value = "1.642.41"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -1.00%  "

$ws.Range("E15").Value = "  -1.79%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.20"
$ws.Range("E16").Value = "  -2.35%  "

$ws.Range("D17").Value = "27.672.43"
$ws.Range("E17").Value = "  +1.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.49"
$ws.Range("E18").Value = "  -0.57%  "

$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.66"
$ws.Range("E20").Value = "  +2.55%  "

$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("E22").Value = "  -1.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.01"
$ws.Range("E23").Value = "  +7.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.96"
$ws.Range("E24").Value = "  -3.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.25"
$ws.Range("E25").Value = "  +1.30%  "

$ws.Range("E26").Value = "  -2.36%  "

$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  -1.46%  "

$ws.Range("E30").Value = "  -0.46%  "

$ws.Range("E31").Value = "  -2.52%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("E33").Value = "  +0.77%  "

$ws.Range("D34").Value = "1.440.78"
$ws.Range("E34").Value = "  +0.10%  "

$ws.Range("E35").Value = "  +1.68%  "

$ws.Range("E36").Value = "  -1.11%  "

$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.885"
$ws.Range("E38").Value = "  -2.75%  "

$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.908"
$ws.Range("E40").Value = "  +15.15%  "

$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.66"
$ws.Range("E43").Value = "  +2.69%  "

$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("E45").Value = "  +2.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.49"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").Value = "1.786.37"
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.09"
$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0106"
$ws.Range("E50").Value = "  -0.38%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0989"
$ws.Range("E51").Value = "  -2.24%  "
